$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count 5 -> 4, Wrong penalty -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right total 70 -> 56, Wrong total -14 -> -28, and fraction label
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -28
$ws.Range("E12").Value = "28 / 112"
